$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

Write-Host ("Overview E col width: " + $overview.Range("E1").ColumnWidth)
Write-Host ("Overview F col width: " + $overview.Range("F1").ColumnWidth)
Write-Host ("zh-cn C col width: " + $zhcn.Range("C1").ColumnWidth)
Write-Host ("de-de C col width: " + $dede.Range("C1").ColumnWidth)
